$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Adjusted capacitor prices a bit" - bump the unit price (col F) for the
# four 100V ceramic capacitor line items (rows 3-6) from $0.002 to $0.30.
# The SUM in B1 and the shared "qty * price" formulas in column G
# recalculate automatically.
$ws.Range("F3").Value = 0.3
$ws.Range("F4").Value = 0.3
$ws.Range("F5").Value = 0.3
$ws.Range("F6").Value = 0.3

# Leave the selection where the author last clicked.
$ws.Range("F20").Select()
